$d = $word.ActiveDocument

# Locate the end of the DOI hyperlink text that should remain as the last
# visible content of the final reference paragraph.
$hyperlinkRange = $d.Content
$foundHyperlink = $hyperlinkRange.Find.Execute(
    "https://doi.org/10.1080/17441692.2021.1912138",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundHyperlink) {
    throw "Could not locate the DOI hyperlink text that anchors the edit."
}
$startPos = $hyperlinkRange.End

# Locate the trailing test sentence that was appended (together with two
# blank manual line breaks before it) after the hyperlink.
$testRange = $d.Content
$foundTest = $testRange.Find.Execute(
    "FAZENDO UM TESTE  DE CRIAR MAIS 3  FRASES",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundTest) {
    throw "Could not locate the trailing test sentence to remove."
}
$endPos = $testRange.End

# Remove everything from right after the hyperlink through the end of the
# test sentence: this deletes the two blank line breaks, the run holding
# the third line break, and the test sentence itself -- restoring the
# paragraph so it ends right after the hyperlink (the paragraph mark
# itself is left untouched).
$deleteRange = $d.Range($startPos, $endPos)
$deleteRange.Delete()
